$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '96.899.17'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.39%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.670.17'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.45%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.28'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.11%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.87'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +10.34%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '654.63'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.58%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.423'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.97%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.08'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.39%  '

# Row 10
$ws.Range("E10").Value = '  +0.06%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.668.80'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.48%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.38'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.93%  '

# Row 14
$ws.Range("E14").Value = '  +6.22%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.356.38'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.52%  '

# Row 16
$ws.Range("E16").Value = '  +3.04%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '96.691.06'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.28%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.90'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.78%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.664.10'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.16%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.85'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.53%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.77'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.32%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.529'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.27%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '533.80'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.48%  '

# Row 24
$ws.Range("E24").Value = '  -0.42%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.22'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +5.17%  '

# Row 26
$ws.Range("E26").Value = '  -0.43%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '102.45'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.34%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.55'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.63%  '

# Row 29
$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.166'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.84%  '

# Row 30
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.39'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +4.25%  '

# Row 31
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.04'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.23%  '

# Row 32
$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.26%  '

# Row 33
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.90'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +14.45%  '

# Row 34
$ws.Range("B34").Value = 'Cronos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.185'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.68%  '

# Row 35
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.40%  '

# Row 36
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '32.68'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +3.08%  '

# Row 37
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '656.95'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +5.74%  '

# Row 38
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.598'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +5.33%  '

# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.89'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.02%  '

# Row 40
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.162'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.72%  '

# Row 41
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.75'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +12.36%  '

# Row 42
$ws.Range("B42").Value = 'ImmutableX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.99'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.28%  '

# Row 43
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.963'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.43%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '38.74'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +17.49%  '

# Row 45
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.04%  '

# Row 46
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0458'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +4.27%  '

# Row 47
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.445'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +8.94%  '

# Row 48
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.33'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.08%  '

# Row 49
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.63'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.09%  '

# Row 50
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.71'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.44%  '

# Row 51
$ws.Range("B51").Value = 'MantraDAO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.66'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +4.50%  '
